$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain text (matches
# the source data's inlineStr representation), so force Text format first.
$textFormatCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.656.12"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "3.896.43"
$ws.Range("E3").Value = "  +6.26%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "608.21"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "173.70"
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("D7").Value = "3.846.05"
$ws.Range("E7").Value = "  +4.97%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  +1.96%  "
$ws.Range("D11").Value = "6.35"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").Value = "0.479"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "39.65"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "0.0000253"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "4.579.80"
$ws.Range("E15").Value = "  +7.05%  "
$ws.Range("D16").Value = "3.934.30"
$ws.Range("E16").Value = "  +7.27%  "
$ws.Range("D17").Value = "69.844.80"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "7.42"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "0.117"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("D20").Value = "16.53"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "503.46"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").Value = "9.63"
$ws.Range("E22").Value = "  +6.45%  "
$ws.Range("D23").Value = "0.748"
$ws.Range("E23").Value = "  +6.20%  "
$ws.Range("D24").Value = "86.66"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  -3.88%  "
$ws.Range("E26").Value = "  +7.87%  "
$ws.Range("D27").Value = "12.60"
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").Value = "10.34"
$ws.Range("E28").Value = "  -7.96%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "2.53"
$ws.Range("E30").Value = "  +4.65%  "
$ws.Range("D31").Value = "2.99"
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("D32").Value = "33.47"
$ws.Range("E32").Value = "  +12.16%  "
$ws.Range("D33").Value = "7.81"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").Value = "0.113"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "6.08"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "1.03"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  +4.54%  "
$ws.Range("D39").Value = "475.29"
$ws.Range("E39").Value = "  +12.56%  "
$ws.Range("D40").Value = "0.332"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "2.04"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "49.66"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").Value = "2.92"
$ws.Range("E43").Value = "  +2.51%  "
$ws.Range("D44").Value = "8.51"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").Value = "42.44"
$ws.Range("E45").Value = "  -5.56%  "
$ws.Range("D46").Value = "2.929.80"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "0.0365"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("D48").Value = "27.35"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").Value = "139.60"
$ws.Range("E49").Value = "  +2.85%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "24.69"
$ws.Range("E51").Value = "  +17.98%  "
